$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 0.0001388888888888889
$ws.Range("K2").Value = 937
$ws.Range("L2").Value = 0.001874

# Row 14 updates
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0.001
$ws.Range("K14").Value = 527
$ws.Range("L14").Value = 0.001054
